$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (Förändrad) holds a date serial value that was bumped by one day
# (45180 -> 45181) for every data row from row 2 through row 390.
$ws.Range("C2:C390").Value = 45181
